# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# style of the existing header (H1 = "IP") and filling in the data
# for rows 2-44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style used by the other header cells (e.g. H1) onto the new
# header cells so they match (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# --- Data rows 2-44 ---
$data = @{
    2  = @(1,3)
    3  = @(1,5)
    4  = @(1,5)
    5  = @(1,5)
    6  = @(1,5)
    7  = @(1,5)
    8  = @(1,5)
    9  = @(1,4)
    10 = @(1,4)
    11 = @(6,8)
    12 = @(7,8)
    13 = @(6,8)
    14 = @(1,3)
    15 = @(8,9)
    16 = @(6,7)
    17 = @(8,9)
    18 = @(4,5)
    19 = @(1,2)
    20 = @(7,7)
    21 = @(6,6)
    22 = @(7,8)
    23 = @(4,5)
    24 = @(6,6)
    25 = @(10,10)
    26 = @(5,6)
    27 = @(8,8)
    28 = @(10,11)
    29 = @(6,7)
    30 = @(5,5)
    31 = @(9,9)
    32 = @(9,9)
    33 = @(9,9)
    34 = @(9,9)
    35 = @(9,9)
    36 = @(8,8)
    37 = @(6,6)
    38 = @(8,8)
    39 = @(8,8)
    40 = @(5,6)
    41 = @(5,6)
    42 = @(6,6)
    43 = @(1,3)
    44 = @(1,2)
}

foreach ($row in 2..44) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}

$wb.Save()
